$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Edit 1: the 'sub' bullet's single run
#   "'sub': pushes the substring of c between b and a"
# is split into four runs and the operands are swapped, so the visible
# text becomes "... between a and b", with "a" and "b" as their own runs:
#   "'sub': pushes the substring of c between " | "a" | " and " | "b"
# -----------------------------------------------------------------------
$rng1 = $d.Content.Duplicate
$found1 = $rng1.Find.Execute("`u{2018}sub`u{2019}: pushes the substring of c between b and a", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Text = ""
    $xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">‘sub’: pushes the substring of c between </w:t></w:r><w:r><w:t>a</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:t>b</w:t></w:r></w:p>
'@
    $rngInsert1 = $d.Range($rng1.Start, $rng1.Start)
    $rngInsert1.InsertXML($xml1)
} else {
    Write-Host "WARNING: could not find the 'sub' bullet text"
}

# -----------------------------------------------------------------------
# Edit 2: two new bulleted list items are added right after the existing
# 'local' bullet (", as local is also accessible from the current
# namespace.") describing the new 'tostack' and 'inverse' stack
# functions.
# -----------------------------------------------------------------------
$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute("as local is also accessible from the current namespace.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>tostack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>’: push a stack containing each individual substring of a to the stack.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>‘inverse’: push an upside down version of a</w:t></w:r></w:p>
'@
    $rngInsert2 = $d.Range($rng2.End, $rng2.End)
    $rngInsert2.InsertXML($xml2)
} else {
    Write-Host "WARNING: could not find the 'local' bullet text"
}

Write-Host "Done"
